# new question and notification function
# Adds "Фамилия" (last name) and "ID" columns, renames "Имя баристы" -> "Имя",
# normalizes the A4 date to a real numeric date, and appends new review rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a literal empty-text value into a cell (not a blank/null
# cell, and not a formula) by using Excel's leading-apostrophe "treat as
# text" prefix, then strip the quote-prefix formatting it implies so the
# cell is left with the default style.
# ---------------------------------------------------------------------------
function Set-EmptyText($row, $col) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'"
    $c.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Header row: rename D1, add E1 "Фамилия" and F1 "ID" with the same header
# style (bold / bordered / centered) as the existing header cells.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 4).Value = "Имя"

$ws.Cells.Item(1, 1).Copy() | Out-Null
$ws.Cells.Item(1, 5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 5).Value = "Фамилия"
$ws.Cells.Item(1, 6).Value = "ID"

# ---------------------------------------------------------------------------
# Existing rows 2-3: add empty E/F cells (typed empty text, matching the
# blank columns that the new header introduced for already-present rows).
# ---------------------------------------------------------------------------
Set-EmptyText 2 5
Set-EmptyText 2 6
Set-EmptyText 3 5
Set-EmptyText 3 6

# ---------------------------------------------------------------------------
# Row 4: the date was previously stored as literal text "2024-05-28"; make
# it a real date value like rows 2-3 (copy their date style), and add the
# new empty E/F cells.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(4, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(4, 1).Value = 45440

Set-EmptyText 4 5
Set-EmptyText 4 6

# ---------------------------------------------------------------------------
# New rows 5-12: full reviews with name/surname/ID columns populated.
# Column A uses the same numeric-date style as rows 2-4.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 5;  Time = "14:00:11.821246"; Point = "Таха Хусейна 2/1"; Name = "Sanzhar"; Surname = "Karibay\"; Id = 507500572 },
    @{ Row = 6;  Time = "14:05:19.681022"; Point = "Таха Хусейна 2/1"; Name = "a";       Surname = "b";        Id = 507500572 },
    @{ Row = 7;  Time = "14:28:39.047375"; Point = "Мухамедханов";     Name = "s";       Surname = "a";        Id = 507500572 },
    @{ Row = 8;  Time = "14:29:56.165677"; Point = "Таха Хусейна 2/1"; Name = "a";       Surname = "a";        Id = 507500572 },
    @{ Row = 9;  Time = "14:30:12.731676"; Point = "Таха Хусейна 2/1"; Name = "a";       Surname = "a";        Id = 507500572 },
    @{ Row = 10; Time = "14:30:43.598356"; Point = "Таха Хусейна 2/1"; Name = "a";       Surname = "a";        Id = 507500572 },
    @{ Row = 11; Time = "14:32:49.106540"; Point = "Таха Хусейна 2/1"; Name = "Sanzhar"; Surname = "Karibay";  Id = 507500572 },
    @{ Row = 12; Time = "14:34:13.195669"; Point = "Таха Хусейна 2/1"; Name = "s";       Surname = "w";        Id = 507500572 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item(2, 1).Copy() | Out-Null
    $ws.Cells.Item($r.Row, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r.Row, 1).Value = 45441

    $ws.Cells.Item($r.Row, 2).Value = $r.Time
    $ws.Cells.Item($r.Row, 3).Value = $r.Point
    $ws.Cells.Item($r.Row, 4).Value = $r.Name
    $ws.Cells.Item($r.Row, 5).Value = $r.Surname
    $ws.Cells.Item($r.Row, 6).Value = $r.Id
}

# ---------------------------------------------------------------------------
# Row 13: another new review, but (like the original row 4) the date arrived
# as plain text rather than a parsed date, so it stays a literal string with
# no special style. A leading apostrophe forces text-entry (otherwise Excel
# auto-parses "2024-05-29" into a real date, same as typing it interactively)
# and Style="Normal" removes the quote-prefix formatting flag afterwards.
# ---------------------------------------------------------------------------
$a13 = $ws.Cells.Item(13, 1)
$a13.Value = "'2024-05-29"
$a13.Style = "Normal"
$ws.Cells.Item(13, 2).Value = "19:35:54.122119"
$ws.Cells.Item(13, 3).Value = "Мухамедханов"
$ws.Cells.Item(13, 4).Value = "sanzhar"
$ws.Cells.Item(13, 5).Value = "karibay"
$ws.Cells.Item(13, 6).Value = 507500572
